# Weekly update of the Jengibre (ginger) price sheet:
# a new market-report row for the week of 2022-10-21 is inserted at row 62,
# pushing the existing rows 62-71 down to 63-72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 62 (shifts rows 62:71 down to 63:72,
# carrying their formatting/styles with them, same as Excel's own
# "Insert Sheet Rows" command).
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with this week's report line.
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44855
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100114007
$ws.Range("G62").Value = "Jengibre"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 500
$ws.Range("K62").Value = 13800
$ws.Range("L62").Value = 14000
$ws.Range("M62").Value = 13900
$ws.Range("N62").Value = "$/caja 13 kilos"
$ws.Range("O62").Value = "Perú"
$ws.Range("P62").Value = 1069
$ws.Range("Q62").Value = 13
$ws.Range("R62").Value = "Hortaliza"
